$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 takes the old values that used to be in row 13
$ws.Range("D2").Value2 = 44585
$ws.Range("M2").Value2 = 160
$ws.Range("N2").Value2 = 6500
$ws.Range("O2").Value2 = 7000
$ws.Range("P2").Value2 = 6750
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value2 = 3375

# Row 3 takes the old values that used to be in row 12
$ws.Range("D3").Value2 = 44627
$ws.Range("M3").Value2 = 45
$ws.Range("N3").Value2 = 6000
$ws.Range("O3").Value2 = 6000
$ws.Range("P3").Value2 = 6000
$ws.Range("R3").Value = "Provincia de Linares"
$ws.Range("S3").Value2 = 3000

# Row 4 takes the old values that used to be in row 8
$ws.Range("D4").Value2 = 44589
$ws.Range("M4").Value2 = 60
$ws.Range("N4").Value2 = 6000
$ws.Range("O4").Value2 = 6000
$ws.Range("P4").Value2 = 6000
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value2 = 3000

# Row 6 takes the old values that used to be in row 4
$ws.Range("D6").Value2 = 44588
$ws.Range("M6").Value2 = 160
$ws.Range("N6").Value2 = 6500
$ws.Range("O6").Value2 = 7000
$ws.Range("P6").Value2 = 6750
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value2 = 3375

# Row 7 takes the old values that used to be in row 2
$ws.Range("D7").Value2 = 44614
$ws.Range("M7").Value2 = 45
$ws.Range("N7").Value2 = 6000
$ws.Range("O7").Value2 = 6000
$ws.Range("P7").Value2 = 6000
$ws.Range("R7").Value = "Provincia de Linares"
$ws.Range("S7").Value2 = 3000

# Row 8 takes the old values that used to be in row 9
$ws.Range("D8").Value2 = 44606
$ws.Range("M8").Value2 = 45
$ws.Range("N8").Value2 = 7000
$ws.Range("O8").Value2 = 7000
$ws.Range("P8").Value2 = 7000
$ws.Range("R8").Value = "Provincia de Linares"
$ws.Range("S8").Value2 = 3500

# Row 9 takes the old values that used to be in row 7
$ws.Range("D9").Value2 = 44586
$ws.Range("M9").Value2 = 80
$ws.Range("N9").Value2 = 7000
$ws.Range("O9").Value2 = 7000
$ws.Range("P9").Value2 = 7000
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value2 = 3500

# Row 10 takes the old values that used to be in row 3
$ws.Range("D10").Value2 = 44582
$ws.Range("M10").Value2 = 150
$ws.Range("N10").Value2 = 6000
$ws.Range("O10").Value2 = 6500
$ws.Range("P10").Value2 = 6233
$ws.Range("R10").Value = "Provincia de Curicó"
$ws.Range("S10").Value2 = 3116

# Row 12 takes the old values that used to be in row 15
$ws.Range("D12").Value2 = 44209
$ws.Range("M12").Value2 = 58
$ws.Range("N12").Value2 = 6000
$ws.Range("O12").Value2 = 6000
$ws.Range("P12").Value2 = 6000
$ws.Range("R12").Value = "Provincia de Curicó"
$ws.Range("S12").Value2 = 3000

# Row 13 takes the old values that used to be in row 14
$ws.Range("D13").Value2 = 44628
$ws.Range("M13").Value2 = 40
$ws.Range("N13").Value2 = 6000
$ws.Range("O13").Value2 = 6000
$ws.Range("P13").Value2 = 6000
$ws.Range("R13").Value = "Provincia de Linares"
$ws.Range("S13").Value2 = 3000

# Row 14 takes the old values that used to be in row 6
$ws.Range("D14").Value2 = 44214
$ws.Range("M14").Value2 = 48
$ws.Range("N14").Value2 = 6000
$ws.Range("O14").Value2 = 6000
$ws.Range("P14").Value2 = 6000
$ws.Range("R14").Value = "Provincia de Linares"
$ws.Range("S14").Value2 = 3000

# Row 15 takes the old values that used to be in row 10
$ws.Range("D15").Value2 = 44592
$ws.Range("M15").Value2 = 30
$ws.Range("N15").Value2 = 8000
$ws.Range("O15").Value2 = 8000
$ws.Range("P15").Value2 = 8000
$ws.Range("R15").Value = "Provincia de Linares"
$ws.Range("S15").Value2 = 4000
